# Apply numeric corrections scraped from the scheduled Kraken_Profits data refresh.
# For each affected sheet/row, update columns H-N (current price & profit calc columns).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 11700.1
$ws.Cells.Item(18, 9).Value = 16917
$ws.Cells.Item(18, 10).Value = 3874.75
$ws.Cells.Item(18, 11).Value = 16917
$ws.Cells.Item(18, 12).Value = 3874.75
$ws.Cells.Item(18, 13).Value = -16633
$ws.Cells.Item(18, 14).Value = -4442.75

$ws.Cells.Item(33, 8).Value = 194
$ws.Cells.Item(33, 9).Value = 194
$ws.Cells.Item(33, 11).Value = 194
$ws.Cells.Item(33, 13).Value = 35

$ws.Cells.Item(40, 8).Value = 4160
$ws.Cells.Item(40, 9).Value = 1100
$ws.Cells.Item(40, 11).Value = 1100
$ws.Cells.Item(40, 13).Value = -925

$ws.Cells.Item(41, 8).Value = 1854.3
$ws.Cells.Item(41, 9).Value = 1942.875
$ws.Cells.Item(41, 10).Value = 1500
$ws.Cells.Item(41, 11).Value = 1942.875
$ws.Cells.Item(41, 12).Value = 1500
$ws.Cells.Item(41, 13).Value = -1502.875
$ws.Cells.Item(41, 14).Value = -2380

$ws.Cells.Item(53, 8).Value = 154.77777
$ws.Cells.Item(53, 10).Value = 226.33333
$ws.Cells.Item(53, 12).Value = 226.33333
$ws.Cells.Item(53, 14).Value = -1500.33333

$ws.Cells.Item(55, 8).Value = 189.63637
$ws.Cells.Item(55, 9).Value = 119
$ws.Cells.Item(55, 11).Value = 119
$ws.Cells.Item(55, 13).Value = 95

$ws.Cells.Item(70, 8).Value = 11277.667
$ws.Cells.Item(70, 9).Value = 1500
$ws.Cells.Item(70, 10).Value = 12499.875
$ws.Cells.Item(70, 11).Value = 4500
$ws.Cells.Item(70, 12).Value = 37499.625
$ws.Cells.Item(70, 13).Value = -4230
$ws.Cells.Item(70, 14).Value = -38039.625

$ws.Cells.Item(73, 8).Value = 11277.667
$ws.Cells.Item(73, 9).Value = 1500
$ws.Cells.Item(73, 10).Value = 12499.875
$ws.Cells.Item(73, 11).Value = 4500
$ws.Cells.Item(73, 12).Value = 37499.625
$ws.Cells.Item(73, 13).Value = -3564
$ws.Cells.Item(73, 14).Value = -39371.625

$ws.Cells.Item(86, 8).Value = 2999.5
$ws.Cells.Item(86, 9).Value = 2999.5
$ws.Cells.Item(86, 11).Value = 2999.5
$ws.Cells.Item(86, 13).Value = -1876.5

$ws.Cells.Item(89, 8).Value = 2999.5
$ws.Cells.Item(89, 9).Value = 2999.5
$ws.Cells.Item(89, 11).Value = 14997.5
$ws.Cells.Item(89, 13).Value = -9381.5

$ws.Cells.Item(94, 8).Value = 888.75
$ws.Cells.Item(94, 9).Value = 888.75
$ws.Cells.Item(94, 11).Value = 888.75
$ws.Cells.Item(94, 13).Value = -437.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 2
$ws.Cells.Item(5, 10).Value = 2
$ws.Cells.Item(5, 12).Value = 2
$ws.Cells.Item(5, 14).Value = -226

$ws.Cells.Item(19, 8).Value = 453.33334
$ws.Cells.Item(19, 9).Value = 453.33334
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 11).Value = 453.33334
$ws.Cells.Item(19, 12).Value = 0
$ws.Cells.Item(19, 13).Value = -224.33334
$ws.Cells.Item(19, 14).Value = ""

$ws.Cells.Item(24, 8).Value = 18371784
$ws.Cells.Item(24, 10).Value = 18371784
$ws.Cells.Item(24, 12).Value = 18371784
$ws.Cells.Item(24, 14).Value = -18372532

$ws.Cells.Item(45, 8).Value = 2183
$ws.Cells.Item(45, 9).Value = 2249.5
$ws.Cells.Item(45, 11).Value = 2249.5
$ws.Cells.Item(45, 13).Value = -1872.5

$ws.Cells.Item(61, 8).Value = 5181.4287
$ws.Cells.Item(61, 9).Value = 4477
$ws.Cells.Item(61, 10).Value = 6942.5
$ws.Cells.Item(61, 11).Value = 4477
$ws.Cells.Item(61, 12).Value = 6942.5
$ws.Cells.Item(61, 13).Value = -4265
$ws.Cells.Item(61, 14).Value = -7366.5

$ws.Cells.Item(74, 8).Value = 5060.4287
$ws.Cells.Item(74, 9).Value = 3004.6
$ws.Cells.Item(74, 10).Value = 10200
$ws.Cells.Item(74, 11).Value = 3004.6
$ws.Cells.Item(74, 12).Value = 10200
$ws.Cells.Item(74, 13).Value = -2130.6
$ws.Cells.Item(74, 14).Value = -11948

$ws.Cells.Item(77, 8).Value = 5060.4287
$ws.Cells.Item(77, 9).Value = 3004.6
$ws.Cells.Item(77, 10).Value = 10200
$ws.Cells.Item(77, 11).Value = 15023
$ws.Cells.Item(77, 12).Value = 51000
$ws.Cells.Item(77, 13).Value = -10655
$ws.Cells.Item(77, 14).Value = -59736

$ws.Cells.Item(100, 8).Value = 18371784
$ws.Cells.Item(100, 10).Value = 18371784
$ws.Cells.Item(100, 12).Value = 18371784
$ws.Cells.Item(100, 14).Value = -18373948

$ws.Cells.Item(102, 8).Value = 2360
$ws.Cells.Item(102, 9).Value = 2040
$ws.Cells.Item(102, 11).Value = 2040
$ws.Cells.Item(102, 13).Value = -418

$ws.Cells.Item(136, 8).Value = 5181.4287
$ws.Cells.Item(136, 9).Value = 4477
$ws.Cells.Item(136, 10).Value = 6942.5
$ws.Cells.Item(136, 11).Value = 13431
$ws.Cells.Item(136, 12).Value = 20827.5
$ws.Cells.Item(136, 13).Value = -10881
$ws.Cells.Item(136, 14).Value = -25927.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 2
$ws.Cells.Item(4, 10).Value = 2
$ws.Cells.Item(4, 12).Value = 2
$ws.Cells.Item(4, 14).Value = -232

$ws.Cells.Item(107, 8).Value = 4000
$ws.Cells.Item(107, 9).Value = 4000
$ws.Cells.Item(107, 11).Value = 4000
$ws.Cells.Item(107, 13).Value = -2080

$ws.Cells.Item(134, 8).Value = 6465.4443
$ws.Cells.Item(134, 9).Value = 5170.2856
$ws.Cells.Item(134, 11).Value = 15510.8568
$ws.Cells.Item(134, 13).Value = -12975.8568

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 10783.625
$ws.Cells.Item(31, 9).Value = 19250
$ws.Cells.Item(31, 10).Value = 7961.5
$ws.Cells.Item(31, 11).Value = 19250
$ws.Cells.Item(31, 12).Value = 7961.5
$ws.Cells.Item(31, 13).Value = -18955
$ws.Cells.Item(31, 14).Value = -8551.5

$ws.Cells.Item(34, 8).Value = 10783.625
$ws.Cells.Item(34, 9).Value = 19250
$ws.Cells.Item(34, 10).Value = 7961.5
$ws.Cells.Item(34, 11).Value = 19250
$ws.Cells.Item(34, 12).Value = 7961.5
$ws.Cells.Item(34, 13).Value = -19048
$ws.Cells.Item(34, 14).Value = -8365.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 667359.4399999999
$ws.Cells.Item(4, 9).Value = 625779.5
$ws.Cells.Item(4, 11).Value = 1877338.5
$ws.Cells.Item(4, 13).Value = -1877226.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 12.25
$ws.Cells.Item(2, 9).Value = 12.25
$ws.Cells.Item(2, 11).Value = 12.25
$ws.Cells.Item(2, 13).Value = 100.75

$ws.Cells.Item(53, 8).Value = 0
$ws.Cells.Item(53, 10).Value = 0
$ws.Cells.Item(53, 12).Value = 0
$ws.Cells.Item(53, 14).Value = ""

$ws.Cells.Item(80, 8).Value = 3070.3333
$ws.Cells.Item(80, 9).Value = 3070.3333
$ws.Cells.Item(80, 11).Value = 3070.3333
$ws.Cells.Item(80, 13).Value = -2072.3333

$ws.Cells.Item(83, 8).Value = 3070.3333
$ws.Cells.Item(83, 9).Value = 3070.3333
$ws.Cells.Item(83, 11).Value = 15351.6665
$ws.Cells.Item(83, 13).Value = -10359.6665

$ws.Cells.Item(107, 8).Value = 332
$ws.Cells.Item(107, 10).Value = 551.5
$ws.Cells.Item(107, 12).Value = 551.5
$ws.Cells.Item(107, 14).Value = -4391.5

$ws.Cells.Item(113, 8).Value = 3011
$ws.Cells.Item(113, 9).Value = 3011
$ws.Cells.Item(113, 11).Value = 3011
$ws.Cells.Item(113, 13).Value = -841

$ws.Cells.Item(122, 8).Value = 944.6667
$ws.Cells.Item(122, 9).Value = 623
$ws.Cells.Item(122, 10).Value = 1266.3334
$ws.Cells.Item(122, 11).Value = 1869
$ws.Cells.Item(122, 12).Value = 3799.0002
$ws.Cells.Item(122, 13).Value = 581
$ws.Cells.Item(122, 14).Value = -8699.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1526.4286
$ws.Cells.Item(16, 9).Value = 1526.4286
$ws.Cells.Item(16, 11).Value = 1526.4286
$ws.Cells.Item(16, 13).Value = -1356.4286

$ws.Cells.Item(46, 8).Value = 5002
$ws.Cells.Item(46, 9).Value = 0
$ws.Cells.Item(46, 10).Value = 5002
$ws.Cells.Item(46, 11).Value = 0
$ws.Cells.Item(46, 12).Value = 5002
$ws.Cells.Item(46, 13).Value = ""
$ws.Cells.Item(46, 14).Value = -5378

$ws.Cells.Item(55, 8).Value = 1154.75
$ws.Cells.Item(55, 9).Value = 986.5714
$ws.Cells.Item(55, 11).Value = 986.5714
$ws.Cells.Item(55, 13).Value = -813.5714

$ws.Cells.Item(74, 8).Value = 20197
$ws.Cells.Item(74, 9).Value = 20197
$ws.Cells.Item(74, 11).Value = 20197
$ws.Cells.Item(74, 13).Value = -19199

$ws.Cells.Item(77, 8).Value = 20197
$ws.Cells.Item(77, 9).Value = 20197
$ws.Cells.Item(77, 11).Value = 60591
$ws.Cells.Item(77, 13).Value = -55599

$ws.Cells.Item(82, 8).Value = 1690.2
$ws.Cells.Item(82, 9).Value = 1690.2
$ws.Cells.Item(82, 11).Value = 1690.2
$ws.Cells.Item(82, 13).Value = -1329.2

$ws.Cells.Item(85, 8).Value = 1690.2
$ws.Cells.Item(85, 9).Value = 1690.2
$ws.Cells.Item(85, 11).Value = 1690.2
$ws.Cells.Item(85, 13).Value = -442.2

$ws.Cells.Item(131, 8).Value = 0
$ws.Cells.Item(131, 10).Value = 0
$ws.Cells.Item(131, 12).Value = 0
$ws.Cells.Item(131, 14).Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(54, 8).Value = 23750
$ws.Cells.Item(54, 9).Value = 17500
$ws.Cells.Item(54, 11).Value = 17500
$ws.Cells.Item(54, 13).Value = -16980

$ws.Cells.Item(62, 8).Value = 2000
$ws.Cells.Item(62, 9).Value = 2000
$ws.Cells.Item(62, 11).Value = 2000
$ws.Cells.Item(62, 13).Value = -1376

$ws.Cells.Item(65, 8).Value = 2000
$ws.Cells.Item(65, 9).Value = 2000
$ws.Cells.Item(65, 11).Value = 10000
$ws.Cells.Item(65, 13).Value = -6880

$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 14).Value = ""

$ws.Cells.Item(136, 8).Value = 2550
$ws.Cells.Item(136, 9).Value = 2550
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 11).Value = 7650
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 13).Value = -5100
$ws.Cells.Item(136, 14).Value = ""
